$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric (non-string) cells first - these don't touch the shared-strings table.
$ws.Range("A20").Value = 23
$ws.Range("B20").Value = 1
$ws.Range("A21").Value = 24
$ws.Range("B21").Value = 4
$ws.Range("A22").Value = 25
$ws.Range("B22").Value = 1
$ws.Range("A23").Value = 26
$ws.Range("B23").Value = 1
$ws.Range("A24").Value = 27
$ws.Range("B24").Value = 10

# String cells written in the precise order needed to reproduce the shared-strings table.
$ws.Range("C20").Value = "TH1"
$ws.Range("F20").Value = "MF52C1103F3380"
$ws.Range("E20").Value = "Thru-hole"
$ws.Range("D20").Value = "NTC"
$ws.Range("G20").Value = "presumably pott'd to batteries?"

$ws.Range("C21").Value = "TP101,TP102, TP103,TP104"
$ws.Range("G21").Value = "wire or keystoone compact TP e.g. p/n 5005"

$ws.Range("C22").Value = "J107"
$ws.Range("D22").Value = "PCB Testing Connector"
$ws.Range("F22").Value = "S2M-110-02-F-D"

$ws.Range("F23").Value = "ITD2-10-D"
$ws.Range("D24").Value = "Mating Connector Pins"
$ws.Range("D23").Value = "Mating Connector Housing"
$ws.Range("F24").Value = "T1M82-L-2426-01-L or T1M82-R-2426-01-L"

# Column width adjustments (columns D and F got wider / best-fit in the real
# edit - this engine's ColumnWidth quantizes to whole pixels, so we feed it
# the input that lands closest to the authored widths of 24.7109375 / 37.42578125).
$ws.Columns("D").ColumnWidth = 23.88
$ws.Columns("F").ColumnWidth = 36.67

# Update the selection to G29 (matches the after-state sheetView selection)
$ws.Range("G29").Select() | Out-Null
